$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Update values in row 2
$ws.Range("A2").Value = 24
$ws.Range("E2").Value = 1627.83
$ws.Range("E2").NumberFormat = "#,##0.00"

# Update values in row 3
$ws.Range("A3").Value = 21
$ws.Range("E3").Value = 8372.17
$ws.Range("E3").NumberFormat = "#,##0.00"
$ws.Range("J3").Value = 8372.17
$ws.Range("J3").NumberFormat = "#,##0.00"

# Update the selected cell on this sheet
$ws.Activate()
$ws.Range("D8").Select()
